# Deploy the implementation guide.
# - Update the "Date" and "Contact" metadata values on the Metadata sheet.
# - Insert a new "Jurisdiction" row (with an empty value) right after
#   "Contact", pushing Description..Count down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update existing metadata values.
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row for "Jurisdiction" right below "Contact" (row 10),
# matching the formatting of the surrounding data rows.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
